# Generate Report for Handback
# Adds a new handback record (463090f6-6b58-48ab-a6c3-9f6fda8abc95.md) as row 3
# to the "Overview", "zh-cn" and "de-de" tables.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Overview" -- add row 3 to table "Overview" (table3.xml)
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.ListRows.Add() | Out-Null

$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c03a0382c7eee7e393a32617fdfdf1863f2f37ce/e2e/463090f6-6b58-48ab-a6c3-9f6fda8abc95.md", "", "", "e2e\463090f6-6b58-48ab-a6c3-9f6fda8abc95.md") | Out-Null

$wsOverview.Range("A3").Value = "463090f6-6b58-48ab-a6c3-9f6fda8abc95.md"
$wsOverview.Range("C3").Value = ".md"
$wsOverview.Range("E3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("G3").Value = "2017-02-17 08:08:03"

# ---------------------------------------------------------------------------
# Sheet "zh-cn" -- add row 3 to table "zh-cn" (table1.xml)
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$loZhCn = $wsZhCn.ListObjects.Item(1)
$loZhCn.ListRows.Add() | Out-Null

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test4-zhcn/blob/fdf7363d1214706a455ba1c266ef5c7face1aaf1/e2e/463090f6-6b58-48ab-a6c3-9f6fda8abc95.md", "", "", "463090f6-6b58-48ab-a6c3-9f6fda8abc95.md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("J3"), "https://github.com/OpenLocalizationTestOrg/ol-test4-zhcn/blob/fdf7363d1214706a455ba1c266ef5c7face1aaf1/e2e/463090f6-6b58-48ab-a6c3-9f6fda8abc95.md", "", "", "463090f6-6b58-48ab-a6c3-9f6fda8abc95.md") | Out-Null

$wsZhCn.Range("B3").Value = ".md"
$wsZhCn.Range("C3").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("D3").Value = "e2e"
$wsZhCn.Range("E3").Value = "ht"
$wsZhCn.Range("F3").Value = "'True"
$wsZhCn.Range("G3").Value = "463090f6-6b58-48ab-a6c3-9f6fda8abc95.f853389c3949ad0ad88f9ca6dd2060afa7864dbb.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2017-02-17 08:07:45"
$wsZhCn.Range("K3").Value = "463090f6-6b58-48ab-a6c3-9f6fda8abc95.f853389c3949ad0ad88f9ca6dd2060afa7864dbb.zh-cn.xlf"
$wsZhCn.Range("L3").Value = "2017-02-17 08:08:39"
$wsZhCn.Range("O3").Value = "'True"
$wsZhCn.Range("Q3").Value = "'False"

# ---------------------------------------------------------------------------
# Sheet "de-de" -- add row 3 to table "de-de" (table2.xml)
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$loDeDe = $wsDeDe.ListObjects.Item(1)
$loDeDe.ListRows.Add() | Out-Null

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test4-dede/blob/b487f93a77a99a5809da56d5058090897c8b675a/e2e/463090f6-6b58-48ab-a6c3-9f6fda8abc95.md", "", "", "463090f6-6b58-48ab-a6c3-9f6fda8abc95.md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("J3"), "https://github.com/OpenLocalizationTestOrg/ol-test4-dede/blob/b487f93a77a99a5809da56d5058090897c8b675a/e2e/463090f6-6b58-48ab-a6c3-9f6fda8abc95.md", "", "", "463090f6-6b58-48ab-a6c3-9f6fda8abc95.md") | Out-Null

$wsDeDe.Range("B3").Value = ".md"
$wsDeDe.Range("C3").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("D3").Value = "e2e"
$wsDeDe.Range("E3").Value = "ht"
$wsDeDe.Range("F3").Value = "'True"
$wsDeDe.Range("G3").Value = "463090f6-6b58-48ab-a6c3-9f6fda8abc95.f853389c3949ad0ad88f9ca6dd2060afa7864dbb.de-de.xlf"
$wsDeDe.Range("H3").Value = "2017-02-17 08:08:03"
$wsDeDe.Range("K3").Value = "463090f6-6b58-48ab-a6c3-9f6fda8abc95.f853389c3949ad0ad88f9ca6dd2060afa7864dbb.de-de.xlf"
$wsDeDe.Range("L3").Value = "2017-02-17 08:09:03"
$wsDeDe.Range("O3").Value = "'True"
$wsDeDe.Range("Q3").Value = "'False"
